# Updated symbol list on Fri Jan 13 21:36:17 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) quotes for the
# coin rows whose figures moved since the previous snapshot. Values are
# written with a leading apostrophe so Excel stores them as text (matching
# the sheet's existing text-formatted numeric/percentage strings) instead of
# auto-converting them to numbers, and the cell style is reset to "Normal"
# afterwards so the quote-prefix flag doesn't linger as a style change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# row -> @{ D = "..."; E = "..." }  (only the columns that actually changed)
$updates = [ordered]@{
    2  = @{ D = "289.76";    E = "0.85%" }
    3  = @{ D = "29.30";     E = "1.79%" }
    4  = @{ D = "5.287";     E = "4.70%" }
    5  = @{ D = "0.07101";   E = "6.55%" }
    6  = @{ D = "7.464";     E = "1.48%" }
    7  = @{ D = "3.568";     E = "5.23%" }
    8  = @{ D = "1.403";     E = "1.88%" }
    9  = @{ D = "0.9081";    E = "-3.52%" }
    10 = @{ D = "0.1618";    E = "4.01%" }
    11 = @{ D = "0.07641";   E = "17.13%" }
    12 = @{ D = "0.07778";   E = "2.42%" }
    13 = @{ D = "0.02918";   E = "-1.12%" }
    14 = @{ D = "0.09036";   E = "0.64%" }
    15 = @{ D = "0.001573";  E = "-2.12%" }
    16 = @{ D = "0.0006532"; E = "1.08%" }
    17 = @{ D = "0.006196";  E = "-2.55%" }
    18 = @{ D = "3.487";     E = "1.17%" }
    19 = @{ D = "2.232";     E = "-1.20%" }
    20 = @{ D = "0.3232";    E = "0.48%" }
    21 = @{ D = "0.1350";    E = "2.93%" }
    22 = @{ D = "4.008";     E = "-1.51%" }
    23 = @{              E = "2.74%" }
    24 = @{ D = "0.04529";   E = "0.77%" }
    25 = @{              E = "2.34%" }
    26 = @{ D = "0.004186";  E = "-6.80%" }
    27 = @{ D = "0.0001168"; E = "-6.51%" }
    28 = @{ D = "0.0001674"; E = "3.28%" }
    40 = @{ D = "0.04389";   E = "4.30%" }
    41 = @{ D = "0.007019";  E = "4.08%" }
    42 = @{ D = "0.1260";    E = "0.37%" }
    43 = @{ D = "0.002067";  E = "2.36%" }
    44 = @{ D = "0.01178";   E = "-4.57%" }
    45 = @{ D = "0.00005842"; E = "1.92%" }
    47 = @{ D = "0.01298";   E = "-0.82%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    if ($cols.Contains("D")) {
        Set-TextValue ("D" + $row) $cols["D"]
    }
    if ($cols.Contains("E")) {
        Set-TextValue ("E" + $row) $cols["E"]
    }
}
